$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1486.5454
$ws.Range("J6").Value = 756.2
$ws.Range("L6").Value = 2268.6
$ws.Range("N6").Value = -2492.6
$ws.Range("H17").Value = 2917.9075
$ws.Range("J17").Value = 2925.327
$ws.Range("L17").Value = 8775.981
$ws.Range("N17").Value = -9111.981
$ws.Range("H80").Value = 50009836
$ws.Range("I80").Value = 125001150
$ws.Range("K80").Value = 375003450
$ws.Range("M80").Value = -375002452
$ws.Range("H83").Value = 50009836
$ws.Range("I83").Value = 125001150
$ws.Range("K83").Value = 1125010350
$ws.Range("M83").Value = -1125005358

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3782.5
$ws.Range("J97").Value = 4630
$ws.Range("L97").Value = 4630
$ws.Range("N97").Value = -5622
$ws.Range("H110").Value = 1639.9231
$ws.Range("I110").Value = 1583.5454
$ws.Range("K110").Value = 1583.5454
$ws.Range("M110").Value = 461.4546
$ws.Range("H122").Value = 2893
$ws.Range("I122").Value = 2102.889
$ws.Range("J122").Value = 3604.1
$ws.Range("K122").Value = 6308.667
$ws.Range("L122").Value = 10812.3
$ws.Range("M122").Value = -3858.667
$ws.Range("N122").Value = -15712.3
$ws.Range("H132").Value = 4635.75
$ws.Range("I132").Value = 4563.5
$ws.Range("K132").Value = 13690.5
$ws.Range("M132").Value = -11160.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1699.3182
$ws.Range("I86").Value = 1568.5
$ws.Range("K86").Value = 1568.5
$ws.Range("M86").Value = -445.5
$ws.Range("H89").Value = 1699.3182
$ws.Range("I89").Value = 1568.5
$ws.Range("K89").Value = 7842.5
$ws.Range("M89").Value = -2226.5
$ws.Range("H94").Value = 3815
$ws.Range("I94").Value = 2271.6667
$ws.Range("J94").Value = 5358.3335
$ws.Range("K94").Value = 2271.6667
$ws.Range("L94").Value = 5358.3335
$ws.Range("M94").Value = -1820.6667
$ws.Range("N94").Value = -6260.3335
$ws.Range("H102").Value = 42880.625
$ws.Range("I102").Value = 11014
$ws.Range("J102").Value = 74747.25
$ws.Range("K102").Value = 11014
$ws.Range("L102").Value = 74747.25
$ws.Range("M102").Value = -7769
$ws.Range("N102").Value = -81237.25
$ws.Range("H134").Value = 2464856
$ws.Range("I134").Value = 2646901
$ws.Range("K134").Value = 7940703
$ws.Range("M134").Value = -7938168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1292.5834
$ws.Range("I16").Value = 1319.1818
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1319.1818
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1032.1818
$ws.Range("N16").Value = -1574
$ws.Range("H31").Value = 5087.3687
$ws.Range("I31").Value = 2119.2222
$ws.Range("K31").Value = 2119.2222
$ws.Range("M31").Value = -1824.2222
$ws.Range("H34").Value = 5087.3687
$ws.Range("I34").Value = 2119.2222
$ws.Range("K34").Value = 2119.2222
$ws.Range("M34").Value = -1917.2222
$ws.Range("H58").Value = 2738.6035
$ws.Range("J58").Value = 4683.8184
$ws.Range("L58").Value = 4683.8184
$ws.Range("N58").Value = -5089.8184
$ws.Range("H60").Value = 6418.2
$ws.Range("J60").Value = 10000
$ws.Range("L60").Value = 10000
$ws.Range("N60").Value = -11022
$ws.Range("H99").Value = 2065
$ws.Range("I99").Value = 2147
$ws.Range("J99").Value = 1901
$ws.Range("K99").Value = 2147
$ws.Range("L99").Value = 1901
$ws.Range("M99").Value = -649
$ws.Range("N99").Value = -4897
$ws.Range("H107").Value = 1811.6666
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 1833.9286
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 1833.9286
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -5673.9286
$ws.Range("H113").Value = 1292.5834
$ws.Range("I113").Value = 1319.1818
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1319.1818
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 850.8181999999999
$ws.Range("N113").Value = -5340
$ws.Range("H126").Value = 2065
$ws.Range("I126").Value = 2147
$ws.Range("J126").Value = 1901
$ws.Range("K126").Value = 6441
$ws.Range("L126").Value = 5703
$ws.Range("M126").Value = -3971
$ws.Range("N126").Value = -10643
$ws.Range("H132").Value = 3805.5
$ws.Range("I132").Value = 3473.7
$ws.Range("K132").Value = 10421.1
$ws.Range("M132").Value = -7891.099999999999
$ws.Range("H134").Value = 1757.1111
$ws.Range("I134").Value = 1767.9395
$ws.Range("J134").Value = 1638
$ws.Range("K134").Value = 5303.818499999999
$ws.Range("L134").Value = 4914
$ws.Range("M134").Value = -2768.818499999999
$ws.Range("N134").Value = -9984
$ws.Range("H136").Value = 2738.6035
$ws.Range("J136").Value = 4683.8184
$ws.Range("L136").Value = 14051.4552
$ws.Range("N136").Value = -19151.4552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 7741.5386
$ws.Range("I2").Value = 34.22222
$ws.Range("J2").Value = 25083
$ws.Range("K2").Value = 34.22222
$ws.Range("L2").Value = 25083
$ws.Range("M2").Value = 78.77778000000001
$ws.Range("N2").Value = -25309
$ws.Range("H102").Value = 3448
$ws.Range("I102").Value = 3369.1428
$ws.Range("K102").Value = 3369.1428
$ws.Range("M102").Value = -1747.1428
$ws.Range("H107").Value = 831.7143
$ws.Range("I107").Value = 817.25
$ws.Range("K107").Value = 817.25
$ws.Range("M107").Value = 1102.75
$ws.Range("H132").Value = 2785.6667
$ws.Range("I132").Value = 2446.375
$ws.Range("K132").Value = 7339.125
$ws.Range("M132").Value = -4809.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3642.2
$ws.Range("I22").Value = 1665.5
$ws.Range("J22").Value = 4960
$ws.Range("K22").Value = 1665.5
$ws.Range("L22").Value = 4960
$ws.Range("M22").Value = -1370.5
$ws.Range("N22").Value = -5550
$ws.Range("H27").Value = 3642.2
$ws.Range("I27").Value = 1665.5
$ws.Range("J27").Value = 4960
$ws.Range("K27").Value = 1665.5
$ws.Range("L27").Value = 4960
$ws.Range("M27").Value = -1558.5
$ws.Range("N27").Value = -5174
$ws.Range("H46").Value = 9487.111000000001
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 10009.883
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 10009.883
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -10385.883
$ws.Range("H61").Value = 2246.9412
$ws.Range("I61").Value = 1913.2
$ws.Range("K61").Value = 1913.2
$ws.Range("M61").Value = -1711.2
$ws.Range("H82").Value = 2780.6191
$ws.Range("I82").Value = 2897.5
$ws.Range("J82").Value = 2624.7778
$ws.Range("K82").Value = 2897.5
$ws.Range("L82").Value = 2624.7778
$ws.Range("M82").Value = -2536.5
$ws.Range("N82").Value = -3346.7778
$ws.Range("H85").Value = 2780.6191
$ws.Range("I85").Value = 2897.5
$ws.Range("J85").Value = 2624.7778
$ws.Range("K85").Value = 2897.5
$ws.Range("L85").Value = 2624.7778
$ws.Range("M85").Value = -1649.5
$ws.Range("N85").Value = -5120.7778
$ws.Range("H113").Value = 2246.9412
$ws.Range("I113").Value = 1913.2
$ws.Range("K113").Value = 1913.2
$ws.Range("M113").Value = 256.8
$ws.Range("H122").Value = 23116.934
$ws.Range("I122").Value = 26985
$ws.Range("J122").Value = 15380.8
$ws.Range("K122").Value = 80955
$ws.Range("L122").Value = 46142.39999999999
$ws.Range("M122").Value = -78505
$ws.Range("N122").Value = -51042.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5046.4
$ws.Range("J81").Value = 5959.8
$ws.Range("L81").Value = 11919.6
$ws.Range("N81").Value = -14041.6
$ws.Range("H84").Value = 5046.4
$ws.Range("J84").Value = 5959.8
$ws.Range("L84").Value = 59598
$ws.Range("N84").Value = -70206
$ws.Range("H130").Value = 59999.5
$ws.Range("J130").Value = 59999.5
$ws.Range("L130").Value = 59999.5
$ws.Range("N130").Value = -70039.5
